$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("O2").Value = 1.28
$ws.Range("Q2").Value = 1.82
$ws.Range("R2").Value = 1.43
$ws.Range("T2").Value = 1.85
$ws.Range("U2").Value = 2.08
$ws.Range("AN2").Value = 9.4

# Row 3 updates
$ws.Range("F3").Value = 1.37
$ws.Range("G3").Value = 1.44
$ws.Range("H3").Value = 10.5
$ws.Range("J3").Value = 4.3
$ws.Range("K3").Value = 5.1
$ws.Range("L3").Value = 1.46
$ws.Range("N3").Value = 2.9
$ws.Range("O3").Value = 1.43
$ws.Range("P3").Value = 1.64
$ws.Range("R3").Value = 1.23
$ws.Range("T3").Value = 2.66
$ws.Range("U3").Value = 1.48
$ws.Range("W3").Value = 3.2
$ws.Range("AC3").Value = 14
$ws.Range("AJ3").Value = 13.5
$ws.Range("AN3").Value = 11

# Row 4 updates
$ws.Range("H4").Value = 1.42
$ws.Range("I4").Value = 1.44
$ws.Range("K4").Value = 5.6
$ws.Range("L4").Value = 1.39
$ws.Range("Q4").Value = 1.85
$ws.Range("S4").Value = 3.15
$ws.Range("V4").Value = 3.25
$ws.Range("X4").Value = 17.5
$ws.Range("Y4").Value = 8.199999999999999
$ws.Range("AH4").Value = 30
$ws.Range("AJ4").Value = 390
$ws.Range("AL4").Value = 160
